# Merging of multiple FC suites into one
#
# The "Input" sheet lists several test orders (TestOrder1..TestOrder5).
# This change fills in the previously-blank Weight/Dimension/Carton values
# for the "TestOrder3" row (row 4), assigns a fresh OrderId to the
# "TestOrder4" row (row 5), and leaves the cursor on A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (TestOrder3): Weight / DimensionL / DimensionW / DimensionH / Cartons
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = 12
$ws.Range("M4").Value = 13
$ws.Range("P4").Value = 1

# Row 5 (TestOrder4): OrderId is reassigned to a new order number.
# Force text formatting first so the numeric-looking id is stored/kept
# as a text value (matching the existing OrderId column's text values)
# rather than being converted to a plain number.
$ws.Range("R5").NumberFormat = "@"
$ws.Range("R5").Value = "51482378"

# Leave the active selection on A4.
$ws.Range("A4").Select()
